# Apply updated TPM-derived NATMI values for Nppc-Npr2.
# Sending clusters now include ECs in addition to FAPs/MuSCs, giving 9 rows
# (3 sending x 3 target clusters) instead of the original 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01599966666666667
$ws.Range("H2").Value = 0.047999
$ws.Range("I2").Value = 0.0009948988194802843
$ws.Range("J2").Value = 0.0009948988194802843
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.65974866666667
$ws.Range("N2").Value = 40.979246
$ws.Range("O2").Value = 0.383515069756994
$ws.Range("P2").Value = 0.3835150697569941
$ws.Range("Q2").Value = 0.2185514254171111
$ws.Range("R2").Value = 1.966962828754
$ws.Range("S2").Value = 0.0003815586901541322
$ws.Range("T2").Value = 0.0003815586901541323

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01599966666666667
$ws.Range("H3").Value = 0.047999
$ws.Range("I3").Value = 0.0009948988194802843
$ws.Range("J3").Value = 0.0009948988194802843
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.104695
$ws.Range("N3").Value = 51.314085
$ws.Range("O3").Value = 0.4802363832729211
$ws.Range("P3").Value = 0.4802363832729211
$ws.Range("Q3").Value = 0.273669418435
$ws.Range("R3").Value = 2.463024765915
$ws.Range("S3").Value = 0.0004777866107897106
$ws.Range("T3").Value = 0.0004777866107897106

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01599966666666667
$ws.Range("H4").Value = 0.047999
$ws.Range("I4").Value = 0.0009948988194802843
$ws.Range("J4").Value = 0.0009948988194802843
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.852797333333333
$ws.Range("N4").Value = 14.558392
$ws.Range("O4").Value = 0.1362485469700849
$ws.Range("P4").Value = 0.1362485469700849
$ws.Range("Q4").Value = 0.07764313973422221
$ws.Range("R4").Value = 0.698788257608
$ws.Range("S4").Value = 0.0001355535185364415
$ws.Range("T4").Value = 0.0001355535185364415

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6973903333333333
$ws.Range("H5").Value = 2.092171
$ws.Range("I5").Value = 0.04336545465636547
$ws.Range("J5").Value = 0.04336545465636547
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.65974866666667
$ws.Range("N5").Value = 40.979246
$ws.Range("O5").Value = 0.383515069756994
$ws.Range("P5").Value = 0.3835150697569941
$ws.Range("Q5").Value = 9.526176675896222
$ws.Range("R5").Value = 85.735590083066
$ws.Range("S5").Value = 0.01663130536757976
$ws.Range("T5").Value = 0.01663130536757976

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nppc"
$ws.Range("C6").Value = "Npr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.6973903333333333
$ws.Range("H6").Value = 2.092171
$ws.Range("I6").Value = 0.04336545465636547
$ws.Range("J6").Value = 0.04336545465636547
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.104695
$ws.Range("N6").Value = 51.314085
$ws.Range("O6").Value = 0.4802363832729211
$ws.Range("P6").Value = 0.4802363832729211
$ws.Range("Q6").Value = 11.928648947615
$ws.Range("R6").Value = 107.357840528535
$ws.Range("S6").Value = 0.02082566910315881
$ws.Range("T6").Value = 0.02082566910315881

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nppc"
$ws.Range("C7").Value = "Npr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.6973903333333333
$ws.Range("H7").Value = 2.092171
$ws.Range("I7").Value = 0.04336545465636547
$ws.Range("J7").Value = 0.04336545465636547
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.852797333333333
$ws.Range("N7").Value = 14.558392
$ws.Range("O7").Value = 0.1362485469700849
$ws.Range("P7").Value = 0.1362485469700849
$ws.Range("Q7").Value = 3.384293949892444
$ws.Range("R7").Value = 30.458645549032
$ws.Range("S7").Value = 0.005908480185626897
$ws.Range("T7").Value = 0.005908480185626897

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Nppc"
$ws.Range("C8").Value = "Npr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.36831233333333
$ws.Range("H8").Value = 46.104937
$ws.Range("I8").Value = 0.9556396465241543
$ws.Range("J8").Value = 0.9556396465241543
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.65974866666667
$ws.Range("N8").Value = 40.979246
$ws.Range("O8").Value = 0.383515069756994
$ws.Range("P8").Value = 0.3835150697569941
$ws.Range("Q8").Value = 209.9272839041669
$ws.Range("R8").Value = 1889.345555137502
$ws.Range("S8").Value = 0.3665022056992602
$ws.Range("T8").Value = 0.3665022056992602

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Nppc"
$ws.Range("C9").Value = "Npr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.36831233333333
$ws.Range("H9").Value = 46.104937
$ws.Range("I9").Value = 0.9556396465241543
$ws.Range("J9").Value = 0.9556396465241543
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.104695
$ws.Range("N9").Value = 51.314085
$ws.Range("O9").Value = 0.4802363832729211
$ws.Range("P9").Value = 0.4802363832729211
$ws.Range("Q9").Value = 262.870295126405
$ws.Range("R9").Value = 2365.832656137645
$ws.Range("S9").Value = 0.4589329275589726
$ws.Range("T9").Value = 0.4589329275589726

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Nppc"
$ws.Range("C10").Value = "Npr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.36831233333333
$ws.Range("H10").Value = 46.104937
$ws.Range("I10").Value = 0.9556396465241543
$ws.Range("J10").Value = 0.9556396465241543
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.852797333333333
$ws.Range("N10").Value = 14.558392
$ws.Range("O10").Value = 0.1362485469700849
$ws.Range("P10").Value = 0.1362485469700849
$ws.Range("Q10").Value = 74.57930510903378
$ws.Range("R10").Value = 671.213745981304
$ws.Range("S10").Value = 0.1302045132659216
$ws.Range("T10").Value = 0.1302045132659216

